$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new rows before current row 10 (the totals row), shifting totals/footer down
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# Step 2: Copy formatting of row 9 (last data row) onto new rows 10 and 11
$ws.Range("A9:N9").Copy()
$ws.Range("A10:N10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A9:N9").Copy()
$ws.Range("A11:N11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Step 3: Fix the price on existing row 9 and its count text
$ws.Range("L9").Value = 11.67
$ws.Range("N9").Value = "0:0"

# Step 4: Fill in new row 10 data (item 7: حفاضات كبار سن ماكسويل 63ق)
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "حفاضات كبار سن ماكسويل 63ق"
$ws.Range("H10").Value = "-1:0"
$ws.Range("L10").Value = 2
$ws.Range("N10").Value = "1:0"
$ws.Rows.Item(10).RowHeight = 25.5

# Step 5: Fill in new row 11 data (item 8: كالونا)
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "كالونا "
$ws.Range("H11").Value = "-1:0"
$ws.Range("L11").Value = 15
$ws.Range("N11").Value = "1:0"
$ws.Rows.Item(11).RowHeight = 24.75

# Step 6: Update totals
$ws.Range("K12").Value = 223.67

# Step 7: Rebuild all merged cell ranges in the canonical order
$ws.Cells.UnMerge()
$ws.Range("C1:L1").Merge()
$ws.Range("E2:F2").Merge()
$ws.Range("G2:I2").Merge()
$ws.Range("J2:L2").Merge()
$ws.Range("B3:G3").Merge()
$ws.Range("H3:K3").Merge()
$ws.Range("L3:M3").Merge()
$ws.Range("B4:G4").Merge()
$ws.Range("H4:K4").Merge()
$ws.Range("L4:M4").Merge()
$ws.Range("B5:G5").Merge()
$ws.Range("H5:K5").Merge()
$ws.Range("L5:M5").Merge()
$ws.Range("B6:G6").Merge()
$ws.Range("H6:K6").Merge()
$ws.Range("L6:M6").Merge()
$ws.Range("B7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("B8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("B9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("B10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("B11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("K12:N12").Merge()
$ws.Range("A13:E13").Merge()
$ws.Range("F13:G13").Merge()
$ws.Range("I13:N13").Merge()

Write-Host "done"
